$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 (index A=2) changes B/C/D values and A value
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = "khuy áo 1111"
$ws.Cells.Item(3, 3).Value = "ka001"
$ws.Cells.Item(3, 4).Value = "K-2-3"

# Row 4
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Mũi khoan  phi 2"
$ws.Cells.Item(4, 3).Value = "2mm"
$ws.Cells.Item(4, 4).Value = "C-02"

# Row 5
$ws.Cells.Item(5, 1).Value = 8
$ws.Cells.Item(5, 2).Value = "Mũi khoan  phi 3"
$ws.Cells.Item(5, 3).Value = "3mm"
$ws.Cells.Item(5, 4).Value = "C-02"

# Row 6
$ws.Cells.Item(6, 1).Value = 9
$ws.Cells.Item(6, 2).Value = "kim khâu"
$ws.Cells.Item(6, 3).Value = "kh002"
$ws.Cells.Item(6, 4).Value = "K-1-2"

# Row 7 (new)
$ws.Cells.Item(7, 1).Value = 15
$ws.Cells.Item(7, 2).Value = "Mũi khoan  phi 5"
$ws.Cells.Item(7, 3).Value = "5mm"
$ws.Cells.Item(7, 4).Value = "C-02"

# Row 8 (new)
$ws.Cells.Item(8, 1).Value = 19
$ws.Cells.Item(8, 2).Value = "Ke vuông"
$ws.Cells.Item(8, 3).Value = "KV120"
$ws.Cells.Item(8, 4).Value = "D-30"

$wb.Save()
